$wb = $excel.ActiveWorkbook

# The "About" sheet (sheet1.xml) gets a new date value in cell C1.
$ws = $wb.Worksheets.Item("About")
$cell = $ws.Range("C1")

# Apply the date number format first so the new style picks up the
# built-in date format (numFmtId 14) instead of auto-registering a
# custom format when the date value is assigned.
$cell.NumberFormat = "mm-dd-yy"

# Serial date 44307 == 2021-04-21
$cell.Value = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0
